# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for the new columns AD, AE, AF
$headers = $ws.Range("AD1:AF1")
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the other header cells (bold, centered/top
# aligned, thin border all around).
$headers.Font.Bold = $true
$headers.HorizontalAlignment = -4108   # xlCenter
$headers.VerticalAlignment = -4160     # xlTop
$headers.Borders.LineStyle = 1         # xlContinuous

# Fill in the team record values for every data row (2 through 60)
$lastRow = 60
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 99   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 63   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
